$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 313, shifting the existing rows 313-326 down to 314-327.
$ws.Rows.Item(313).Insert()

# Populate the newly inserted row 313 with a new weekly record.
$ws.Range("A313").Value = 3
$ws.Range("B313").Value = "Femacal de La Calera"
$ws.Range("C313").Value = "Coquimbo"
$ws.Range("D313").Value = 44509
$ws.Range("E313").Value = 5
$ws.Range("F313").Value = 100112037
$ws.Range("G313").Value = "Cebollín"
$ws.Range("H313").Value = "Sin especificar"
$ws.Range("I313").Value = "Primera"
$ws.Range("J313").Value = 160
$ws.Range("K313").Value = 3000
$ws.Range("L313").Value = 3000
$ws.Range("M313").Value = 3000
$ws.Range("N313").Value = "$/paquete 36 unidades"
$ws.Range("O313").Value = "Provincia de Quillota"
$ws.Range("P313").Value = 83
$ws.Range("Q313").Value = 36
$ws.Range("R313").Value = "Hortaliza"
